$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$excel.ActiveWindow.Zoom = 140
$ws.Range("B2").Select()
